$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 562, pushing the existing row 562 (and
# everything below it) down by one. This grows the used range from
# A1:R684 to A1:R685.
$ws.Rows(562).Insert()

# Populate the newly inserted row 562 with the new weekly record.
$ws.Range("A562").Value = 9
$ws.Range("B562").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C562").Value = "Metropolitana"
$ws.Range("D562").Value = 44889
$ws.Range("E562").Value = 13
$ws.Range("F562").Value = 100114014
$ws.Range("G562").Value = "Betarraga"
$ws.Range("H562").Value = "Sin especificar"
$ws.Range("I562").Value = "Primera"
$ws.Range("J562").Value = 15000
$ws.Range("K562").Value = 80
$ws.Range("L562").Value = 100
$ws.Range("M562").Value = 88
$ws.Range("N562").Value = "$/unidad"
$ws.Range("O562").Value = "Región Metropolitana"
$ws.Range("P562").Value = 88
$ws.Range("Q562").Value = 1
$ws.Range("R562").Value = "Hortaliza"
